# Regression_TC.xlsx - "adding api and access"
#  - Unhide / widen the "Active" (AB) column so it's visible again.
#  - Sync the Active (AB) flags with the actual Status (AC) text for each
#    test case: TRUE (green "Good" style) when Status is complete,
#    FALSE (amber "Neutral" style) when Status still needs attention.
#  - Leave the selection on the last touched cell (AB16).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Regression")

# --- Unhide / widen column AB ("Active") --------------------------------
$ws.Columns.Item(28).Hidden = $false
$ws.Columns.Item(28).ColumnWidth = 18.71

# --- Rows whose Active flag simply flips to TRUE (style already "Good") -
$trueRows = 2, 3, 4, 5, 6
foreach ($r in $trueRows) {
    $ws.Cells.Item($r, 28).Value = $true
}

# --- Rows that become TRUE / "Good" (green) style ------------------------
# Pull the known-good green formatting from AB2 (already styled "Good").
$greenRows = 8, 9, 10, 11, 12, 15
$ws.Range("AB2").Copy() | Out-Null
foreach ($r in $greenRows) {
    $cell = $ws.Cells.Item($r, 28)
    $cell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}
$ws.Cells.Item(15, 28).Value = $true

# --- Rows that become FALSE / "Neutral" (amber) style ---------------------
# Pull the known amber formatting from AC7 (already styled "Neutral").
$amberRows = 7, 13, 14, 16
$ws.Range("AC7").Copy() | Out-Null
foreach ($r in $amberRows) {
    $cell = $ws.Cells.Item($r, 28)
    $cell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}
$ws.Cells.Item(13, 28).Value = $false
$ws.Cells.Item(16, 28).Value = $false

$excel.CutCopyMode = $false

# --- Leave the cursor where the author left it ---------------------------
$ws.Range("AB16").Select() | Out-Null
